# Insert a new weekly price record as row 19, pushing existing rows 19-95
# down to 20-96 (dimension grows from A1:T95 to A1:T96).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 19:95 down by one row (standard Excel "insert row" behaviour,
# which also carries the D-column date style down with the data).
$ws.Rows(19).Insert()

# Populate the newly-inserted row 19 with the new weekly record.
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value = 44453
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100108
$ws.Range("H19").Value = "Tropicales y subtropicales"
$ws.Range("I19").Value = 100108005
$ws.Range("J19").Value = "Piña"
$ws.Range("K19").Value = "Caramelo"
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 21000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 21500
$ws.Range("Q19").Value = "$/caja 14 unidades"
$ws.Range("R19").Value = "Ecuador"
$ws.Range("S19").Value = 1536
$ws.Range("T19").Value = 14
